$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('3535 Opal Meadow Heights Aged Care Community Meadow Heights', 28),
    @('ABD Group 512 Melbourne Road Construction Site Spotswood', 5),
    @('Acquire BPO Southbank', 7),
    @('Adorn Cosmetics Clayton', 5),
    @('Al Haj Halal Meats Glenroy', 34),
    @('Al-Taqwa College Truganina', 5),
    @('Best&Less Fountain Gate Narre Warren', 5),
    @('Budget Car and Truck Rentals Campbellfield', 7),
    @('CS Square Caroline Springs', 9),
    @('Cannie Road Construction Site Cannie', 7),
    @('Caroline Springs Police Station', 12),
    @('Cedars Medical Clinic Coburg', 41),
    @('Cedars Medical Clinic Coburg', 28),
    @('Chemist Warehouse Campbellfield DC', 5),
    @('Chemist Warehouse Fillo Drive Somerton', 5),
    @('City of Moreland Community', 6),
    @('City of Wyndham Community', 6),
    @('Classy Cabinets and Kitchens Craigieburn', 10),
    @('Coles Aurora Village Epping', 6),
    @('Coles Broadmeadows Central Shopping Centre', 9),
    @('Coles Campbellfield Plaza Campbellfield', 9),
    @('Coles Coburg North Village', 26),
    @('Coles Greenvale Shopping Centre', 7),
    @('Coles Pakenham Place Shopping Centre', 11),
    @('Coles Roxburgh Village Roxburgh Park', 11),
    @('Community Kids Meadow Heights', 12),
    @('Construction Site Olea Apartment Caulfield North', 16),
    @('Costco Wholesale Epping', 24),
    @('Crusader Caravans Epping', 14),
    @('DRC Laverton Automotive Repairs Laverton North', 5),
    @('Direct Freight Express Cambellfield', 13),
    @('Epworth Healthcare Epworth Richmond Emergency Department', 5),
    @('Fitzroy Community School Fitzroy North', 36),
    @('Fonterra Manufacturing Workplace Campbellfield', 7),
    @('Glenroy West Primary School', 6),
    @('Goodstart Early Learning Altona', 5),
    @('Green Leaves Early Learning Centre Highlands Craigieburn', 7),
    @('Gumboots Early Learning Centre South Morang', 5),
    @('Hamilton Marino 236 Jasper Road McKinnon', 12),
    @('Health Care Providers Association South Melbourne', 13),
    @('IGA Meadow Heights Shopping Centre Meadow Heights', 6),
    @('ISS Factory Level 1 Terminal 2 Melbourne Airport Tullamarine', 9),
    @('Ibis Kingsgate Hotel Melbourne', 6),
    @('Ilim College Glenroy Campus Hadfield', 16),
    @('Ilim Learning Sanctuary Glenroy', 11),
    @('Impact Designer Homes Epping', 5),
    @('Industrial Galvanizers Valmont Coatings Campbellfield', 22),
    @('Islamic College of Melbourne Tarneit', 5),
    @('KFC Fawkner', 8),
    @('Kasr Sweets Coolaroo', 5),
    @('Kids House Early Learning Cheltenham', 12),
    @('Learning Nest Early Learning Centre Meadow Heights', 5),
    @('Level Crossing Removal Project Lilydale Construction Site John Street', 9),
    @('Lineage Logistics Laverton North', 8),
    @('Linfox Somerton National Distribution Centre Somerton', 9),
    @('McDonalds Thomastown II', 7),
    @('Melbourne Metropolitan Remand Centre Ravenhall', 11),
    @('Melbourne Truck Repairs Campbellfield', 7),
    @('Melbourne West Police Station Docklands', 7),
    @('Melbourne Youth Justice Centre Parkville', 5),
    @('Melton Police Station Melton', 5),
    @('Mercy Hospital for Women Heidelberg', 5),
    @('Mernda YMCA Early Learning Centre Mernda', 5),
    @('Montessori Beginnings Greenvale', 5),
    @('MyCentre Childcare Broadmeadows', 17),
    @('National Gallery of Victoria Melbourne', 9),
    @('Newbury Child and Community Centre Craigieburn', 5),
    @('Nido Early School Moonee Ponds', 14),
    @('North Geelong House Party', 7),
    @('Northern Health Northern Hospital Epping Emergency Department Tier 1B', 47),
    @('Northern Health The Northern Hospital Epping', 5),
    @('OnQ Plumbing and Excavations Craigieburn', 18),
    @('Oporto Coolaroo', 13),
    @('Our Lady Help of Christian''s Primary School Brunswick East', 8),
    @('Paisley Park Early Learning Centre Bundoora', 8),
    @('Panorama Construction Site Whitehorse Rd Box Hill', 17),
    @('Private Residence Northern Community Services Fawkner', 5),
    @('Ramsay Health Care Warringal Private Hospital Heidelberg', 9),
    @('Richmond Quarter 261-271 Bridge Road Construction Site Richmond', 10),
    @('Salta Drive Construction Site Rangedale Drainage Altona North', 7),
    @('Sharpline Stainless Steel Coburg North', 5),
    @('St Vincents Hospital Emergency Department Melbourne', 5),
    @('Tek Foods Somerton', 12),
    @('The Homestead Child and Family Centre Roxburgh Park', 12),
    @('The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B', 11),
    @('ThorwestenCabinets Pakenham', 13),
    @('Total Window Concepts Hoppers Crossing', 6),
    @('Unilodge College Square Student Accommodation 570 Lygon Street Carlton', 14),
    @('Werribee Mercy Hospital Emergency Department', 8),
    @('Western Health Footscray Hospital Emergency Department', 6),
    @('Western Health Sunshine Hospital Emergency Department', 8),
    @('Woodlands Long Day Care and Kindergarten Roxburgh Park', 5),
    @('Woolworths Greenvale Lakes Roxburgh Park', 5),
    @('Yara Childcare Centre Truganina', 10)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

Write-Output "done"